$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds date serials, columns J/M/P hold numeric volume/price data.
# Update each cell directly with its new value per the target diff.

$ws.Range("D2").Value  = 44188
$ws.Range("J2").Value  = 200

$ws.Range("D3").Value  = 44188
$ws.Range("J3").Value  = 100

$ws.Range("D4").Value  = 44230
$ws.Range("D5").Value  = 44230

$ws.Range("D6").Value  = 44335
$ws.Range("J6").Value  = 150
$ws.Range("M6").Value  = 633
$ws.Range("P6").Value  = 106

$ws.Range("D7").Value  = 44335
$ws.Range("J7").Value  = 50

$ws.Range("D8").Value  = 44358
$ws.Range("J8").Value  = 200
$ws.Range("M8").Value  = 650
$ws.Range("P8").Value  = 108

$ws.Range("D9").Value  = 44358
$ws.Range("J9").Value  = 100

$ws.Range("D10").Value = 44328
$ws.Range("J10").Value = 100

$ws.Range("D11").Value = 44328
$ws.Range("J11").Value = 50

$ws.Range("D12").Value = 44293
$ws.Range("J12").Value = 100

$ws.Range("D13").Value = 44293
$ws.Range("J13").Value = 50

$ws.Range("D14").Value = 44308
$ws.Range("J14").Value = 200

$ws.Range("D15").Value = 44308
$ws.Range("J15").Value = 100

$ws.Range("D16").Value = 44321
$ws.Range("D17").Value = 44321
